# "Provided some page load checks"
#
# Adds a new "Bill" table definition block in column U (mirroring the
# existing Staff/ServiceRequest/OnlineBooking/Booking table-definition
# blocks already present in columns G/J/M/N), and extends the Room
# "Status" enum legend (column F, rows 21-24) with a new "Booked" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New enum value for the Room.Status legend (column F) ---------------
# F21 = "Status" header, F22:F24 = UnderMaintenance/Empty/Occupied (existing,
# italic style). Add F25 = "Booked" using the same italic style as F22:F24.
# Written before the "Bill" strings below so the new shared-string table
# keeps the same ordering as produced by the original edit (Booked first).
$ws.Range("F25").Value = "Booked"
$ws.Range("F25").Font.Italic = $true

# --- New "Bill" table definition (column U, rows 1-10) -------------------
$ws.Range("U1").Value = "Bill"
$ws.Range("U2").Value = "BillID"
$ws.Range("U3").Value = "BookingID"
$ws.Range("U4").Value = "RoomCharges"
$ws.Range("U5").Value = "ServiceCharges"
$ws.Range("U6").Value = "Total"
$ws.Range("U7").Value = "DiscountedAmount"
$ws.Range("U8").Value = "PaybleAmount"
$ws.Range("U9").Value = "PaymentMode"
$ws.Range("U10").Value = "Details"

# --- Column sizing --------------------------------------------------------
# Column F widens to fit the longer "UnderMaintenance"/"Booked" labels, and
# the new column U gets sized to fit its longest label ("DiscountedAmount").
$ws.Columns.Item(6).ColumnWidth = 18
$ws.Columns.Item(21).ColumnWidth = 14.166666666666666

# --- View state ------------------------------------------------------------
# Select the BookingID cell of the new Bill block, matching where the
# editor's cursor ended up after authoring the new table.
$ws.Range("U3").Select()
